# EscolherConfiguracaoOtima.xlsx - "Primeira versão dos use cases done"
# Fills in the use-case template (pre/post conditions, normal scenario steps,
# alternative scenario and exception) and appends a new "Exceção" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Resize the body font from 14pt to 11pt (everything except the two
#    merged "scenario label" cells in column B for rows 17-19, which keep
#    their original font until the new block is copied from them below).
# ---------------------------------------------------------------------------
$ws.Range("B2:D16").Font.Size = 11
$ws.Range("C17:D19").Font.Size = 11

# ---------------------------------------------------------------------------
# 2) Fill in the missing use-case content.
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = "Todos os componentes primários encontram-se selecionados"

$ws.Range("C7").Value = "1. Insere valor que está disposto a gastar"
$ws.Range("D8").Value = "2. Calcula configuração que maximiza o valor a gastar dentro do limite indicado"
$ws.Range("D9").Value = "3. Apresenta configuração sugerida"
$ws.Range("C10").Value = "4. Aceita a configuração proposta"
$ws.Range("D11").Value = "5. Adiciona os componentes aos componentes selecionados"

$ws.Range("B17").Value = "Cenário Alternativo 1 [não aceita configuração proposta] (passo 3)"
$ws.Range("D17").Value = "3.1 Requisita a inserção de um novo valor"
$ws.Range("D18").Value = "3.2 Regressa a 1"

# ---------------------------------------------------------------------------
# 3) Wrap text for the cells whose content needs more than one line, and
#    resize the corresponding rows to fit.
# ---------------------------------------------------------------------------
$ws.Range("C7").WrapText = $true
$ws.Range("D8").WrapText = $true
$ws.Range("D11").WrapText = $true
$ws.Range("B17:B19").WrapText = $true
$ws.Range("D17").WrapText = $true

$ws.Rows(7).RowHeight = 38.25
$ws.Rows(8).RowHeight = 57
$ws.Rows(11).RowHeight = 38.25
$ws.Rows(17).RowHeight = 38.25

# ---------------------------------------------------------------------------
# 4) Add the new "Exceção 2" block (rows 20-22) by duplicating the format of
#    the "Cenário Alternativo 1" block (rows 17-19), then filling its text.
#    Merge the destination cells first so PasteSpecial keeps a uniform
#    (un-split) border across the whole merged area.
# ---------------------------------------------------------------------------
$ws.Range("B20:B22").Merge()
$ws.Range("B17:D19").Copy()
$ws.Range("B20").PasteSpecial(-4122)

$ws.Range("B20:B22").Font.Size = 12

$ws.Range("B20").Value = "Exceção 2 [não insere novo valor] (passo 3.1)"
$ws.Range("D20").Value = "3.1.1 Informa de insucesso a calcular configuração ótima"

$ws.Rows(20).RowHeight = 38.25

# ---------------------------------------------------------------------------
# 5) Page setup + selection.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("C20").Select()
